$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H31").Value = 793.3333
$ws.Range("I31").Value = 793.3333
$ws.Range("K31").Value = 2379.9999
$ws.Range("M31").Value = -2149.9999

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H39").Value = 315.94736
$ws.Range("I39").Value = 150.57143
$ws.Range("J39").Value = 412.41666
$ws.Range("K39").Value = 451.71429
$ws.Range("L39").Value = 1237.24998
$ws.Range("M39").Value = -155.71429
$ws.Range("N39").Value = -1829.24998

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 2830.8462
$ws.Range("I40").Value = 4800
$ws.Range("J40").Value = 1600.125
$ws.Range("K40").Value = 4800
$ws.Range("L40").Value = 1600.125
$ws.Range("M40").Value = -4625
$ws.Range("N40").Value = -1950.125

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H98").Value = 5445.8184
$ws.Range("I98").Value = 5211.5557
$ws.Range("K98").Value = 5211.5557
$ws.Range("M98").Value = -3713.5557

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H103").Value = 675
$ws.Range("J103").Value = 675
$ws.Range("L103").Value = 2025
$ws.Range("N103").Value = -3197

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H122").Value = 5445.8184
$ws.Range("I122").Value = 5211.5557
$ws.Range("K122").Value = 15634.6671
$ws.Range("M122").Value = -13184.6671

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 619726.7
$ws.Range("I137").Value = 2046.75
$ws.Range("J137").Value = 1284920.5
$ws.Range("K137").Value = 6140.25
$ws.Range("L137").Value = 3854761.5
$ws.Range("M137").Value = -3590.25
$ws.Range("N137").Value = -3859861.5

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 3455.7646
$ws.Range("I138").Value = 2475.9092
$ws.Range("J138").Value = 3924.3914
$ws.Range("K138").Value = 7427.7276
$ws.Range("L138").Value = 11773.1742
$ws.Range("M138").Value = -2287.7276
$ws.Range("N138").Value = -22053.1742

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 4845.7837
$ws.Range("I61").Value = 5152.0435
$ws.Range("K61").Value = 5152.0435
$ws.Range("M61").Value = -4940.0435

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H130").Value = 56601.875
$ws.Range("J130").Value = 56601.875
$ws.Range("L130").Value = 56601.875
$ws.Range("N130").Value = -66641.875

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 2698.147
$ws.Range("I132").Value = 2371.5417
$ws.Range("K132").Value = 7114.625100000001
$ws.Range("M132").Value = -4584.625100000001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 4845.7837
$ws.Range("I136").Value = 5152.0435
$ws.Range("K136").Value = 15456.1305
$ws.Range("M136").Value = -12906.1305

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 980236.5
$ws.Range("I105").Value = 1253574.4
$ws.Range("J105").Value = 4029.7144
$ws.Range("K105").Value = 1253574.4
$ws.Range("L105").Value = 4029.7144
$ws.Range("M105").Value = -1251827.4
$ws.Range("N105").Value = -7523.7144

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H115").Value = 10000
$ws.Range("J115").Value = 10000
$ws.Range("L115").Value = 10000
$ws.Range("N115").Value = -13134

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 2209.9092
$ws.Range("I134").Value = 2800.125
$ws.Range("J134").Value = 636
$ws.Range("K134").Value = 8400.375
$ws.Range("L134").Value = 1908
$ws.Range("M134").Value = -5865.375
$ws.Range("N134").Value = -6978

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H4").Value = 35862.07

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 689463.9399999999
$ws.Range("I31").Value = 5924.25
$ws.Range("J31").Value = 1091546.1
$ws.Range("K31").Value = 5924.25
$ws.Range("L31").Value = 1091546.1
$ws.Range("M31").Value = -5629.25
$ws.Range("N31").Value = -1092136.1

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 689463.9399999999
$ws.Range("I34").Value = 5924.25
$ws.Range("J34").Value = 1091546.1
$ws.Range("K34").Value = 5924.25
$ws.Range("L34").Value = 1091546.1
$ws.Range("M34").Value = -5722.25
$ws.Range("N34").Value = -1091950.1

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H94").Value = 237.14285
$ws.Range("I94").Value = 77
$ws.Range("K94").Value = 77
$ws.Range("M94").Value = 374

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H122").Value = 4779.1665
$ws.Range("I122").Value = 4639.3335
$ws.Range("J122").Value = 4919
$ws.Range("K122").Value = 13918.0005
$ws.Range("L122").Value = 14757
$ws.Range("M122").Value = -11468.0005
$ws.Range("N122").Value = -19657

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 39.29032
$ws.Range("I2").Value = 32.333332
$ws.Range("J2").Value = 40.96
$ws.Range("K2").Value = 193.999992
$ws.Range("L2").Value = 245.76
$ws.Range("M2").Value = -80.99999199999999
$ws.Range("N2").Value = -471.76

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H14").Value = 1421.8
$ws.Range("I14").Value = 1421.8
$ws.Range("K14").Value = 4265.4
$ws.Range("M14").Value = -4092.4

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H33").Value = 136.8
$ws.Range("I33").Value = 0
$ws.Range("J33").Value = 136.8
$ws.Range("K33").Value = 0
$ws.Range("L33").Value = 820.8000000000001
$ws.Range("M33").ClearContents()
$ws.Range("N33").Value = -1386.8

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 2432.6235
$ws.Range("J68").Value = 4327.3716
$ws.Range("L68").Value = 12982.1148
$ws.Range("N68").Value = -14604.1148

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H71").Value = 2432.6235
$ws.Range("J71").Value = 4327.3716
$ws.Range("L71").Value = 38946.3444
$ws.Range("N71").Value = -47058.3444

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H92").Value = 260.77777
$ws.Range("I92").Value = 241
$ws.Range("J92").Value = 266.42856
$ws.Range("K92").Value = 723
$ws.Range("L92").Value = 799.28568
$ws.Range("M92").Value = 525
$ws.Range("N92").Value = -3295.28568

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H95").Value = 6006
$ws.Range("J95").Value = 6000
$ws.Range("L95").Value = 18000
$ws.Range("N95").Value = -22118

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H107").Value = 390.52112
$ws.Range("I107").Value = 262.29688
$ws.Range("J107").Value = 1562.8572
$ws.Range("K107").Value = 786.89064
$ws.Range("L107").Value = 4688.571599999999
$ws.Range("M107").Value = 1133.10936
$ws.Range("N107").Value = -8528.571599999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 1132.1562
$ws.Range("J131").Value = 955.2941
$ws.Range("L131").Value = 2865.8823
$ws.Range("N131").Value = -12945.8823

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H139").Value = 3916609.2
$ws.Range("I139").Value = 11741750
$ws.Range("J139").Value = 4039.1667
$ws.Range("K139").Value = 35225250
$ws.Range("L139").Value = 12117.5001
$ws.Range("M139").Value = -35220110
$ws.Range("N139").Value = -22397.5001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H35").Value = 12500
$ws.Range("J35").Value = 12500
$ws.Range("L35").Value = 12500
$ws.Range("N35").Value = -13096

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 6090.5713
$ws.Range("I132").Value = 5905.5
$ws.Range("J132").Value = 6337.3335
$ws.Range("K132").Value = 17716.5
$ws.Range("L132").Value = 19012.0005
$ws.Range("M132").Value = -15186.5
$ws.Range("N132").Value = -24072.0005

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 3221.0356
$ws.Range("I40").Value = 3040.9
$ws.Range("J40").Value = 3671.375
$ws.Range("K40").Value = 3040.9
$ws.Range("L40").Value = 3671.375
$ws.Range("M40").Value = -2904.9
$ws.Range("N40").Value = -3943.375

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H128").Value = 57900
$ws.Range("J128").Value = 57900
$ws.Range("L128").Value = 57900
$ws.Range("N128").Value = -67860

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H2").Value = 8918667
$ws.Range("I2").Value = 35125000
$ws.Range("J2").Value = 1431142.9
$ws.Range("K2").Value = 35125000
$ws.Range("L2").Value = 1431142.9
$ws.Range("M2").Value = -35124888
$ws.Range("N2").Value = -1431366.9

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H124").Value = 59200
$ws.Range("J124").Value = 59200
$ws.Range("L124").Value = 59200
$ws.Range("N124").Value = -69020

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 4617.5
$ws.Range("J136").Value = 7389.593
$ws.Range("L136").Value = 22168.779
$ws.Range("N136").Value = -27268.779
